$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Delete the "Infant testing" (row 30) and "Mother testing" (row 31) rows,
# which were erroneous/duplicate entries — everything below shifts up.
$ws.Rows.Item(30).Resize(2).Delete()
